$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected; unprotect to make the edits, then re-protect
# so the resulting workbook still reports itself as protected.
$ws.Unprotect()

$ws.Range("A13").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-13 for illustrative purposes only and are subject to change."

$ws.Range("D2").Value = 0.09051529970864942
$ws.Range("E2").Value = -0.01800240032004274
$ws.Range("D3").Value = 0.1058772418873681
$ws.Range("E3").Value = 0.002943202454755989
$ws.Range("D4").Value = 0.1205952072066437
$ws.Range("E4").Value = 0.0140618722378465
$ws.Range("D5").Value = 0.1415103583694842
$ws.Range("E5").Value = 0.01282078512808527
$ws.Range("D6").Value = 0.1390627906916728
$ws.Range("E6").Value = 0.01279347673274289
$ws.Range("D7").Value = 0.1471449722848034
$ws.Range("E7").Value = 0.01698369565217384
$ws.Range("D8").Value = 0.1253969829180922
$ws.Range("E8").Value = 0.01153726223885254
$ws.Range("D9").Value = 0.1298971469332861
$ws.Range("E9").Value = 0.009251821855608355
$ws.Range("D10").Value = 0.9999999999999999
$ws.Range("E10").Value = 0.009118878934074148

$ws.Protect()
